$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 23 and 24 (budget = 900000000 entries no longer solved/printed)
$ws.Rows("23:24").Delete()

# Update Solving Time (s) column F for remaining rows 3-22 with new values
$newSolvingTimes = @(
    10.62743169999885,
    11.21693030000461,
    10.892906699999,
    11.67338169999857,
    10.88322079999489,
    10.96941220000008,
    11.10764560000098,
    11.28365040000062,
    11.3287593999994,
    5.835396100002981,
    5.152236800000537,
    5.153378800001519,
    4.901795699996001,
    4.952405299998645,
    5.099446799998987,
    5.031342300004326,
    4.992929000000004,
    4.920939700001327,
    5.113470500000403,
    4.999486900000193
)

for ($i = 0; $i -lt $newSolvingTimes.Count; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 6).Value = $newSolvingTimes[$i]
}
